$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C3").Value = -12.1229
$ws.Range("A9").Value = -20.47729999999998
$ws.Range("C11").Value = -13.7345
$ws.Range("A18").Value = -22.89990000000001
$ws.Range("A20").Value = -22.23070000000003
$ws.Range("E21").Value = 13.05169999999999
